$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 26
$ws.Range("F3").Value = 1006
$ws.Range("G3").Value = 79
$ws.Range("F5").Value = 1173
$ws.Range("F6").Value = 966
$ws.Range("F7").Value = 295
$ws.Range("F10").Value = 913
$ws.Range("F11").Value = 333
$ws.Range("F13").Value = 535
$ws.Range("F14").Value = 1387
$ws.Range("F16").Value = 1294
$ws.Range("F17").Value = 2951
$ws.Range("F18").Value = 348
$ws.Range("F19").Value = 1575
$ws.Range("F20").Value = 1325
$ws.Range("F22").Value = 221
$ws.Range("I23").Value = "//i2.hdslb.com/bfs/openplatform/202404/A9SL4gE21713522395699.png"
$ws.Range("F24").Value = 244
$ws.Range("F26").Value = 1084
$ws.Range("F27").Value = 379
$ws.Range("F28").Value = 3365
$ws.Range("F29").Value = 654
$ws.Range("F31").Value = 1490

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 166
$ws.Range("F8").Value = 13
$ws.Range("F10").Value = 40
$ws.Range("F11").Value = 11
$ws.Range("F13").Value = 4
$ws.Range("F15").Value = 43
$ws.Range("F17").Value = 1

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 26
$ws.Range("F6").Value = 1006
$ws.Range("G6").Value = 79
$ws.Range("F9").Value = 1173
$ws.Range("F10").Value = 966
$ws.Range("F11").Value = 295
$ws.Range("F16").Value = 166
$ws.Range("F17").Value = 13
$ws.Range("F19").Value = 40
$ws.Range("F21").Value = 11
$ws.Range("F23").Value = 913
$ws.Range("F24").Value = 333
$ws.Range("F26").Value = 535
$ws.Range("F27").Value = 1387
$ws.Range("F29").Value = 1294
$ws.Range("F30").Value = 2951
$ws.Range("F31").Value = 348
$ws.Range("F32").Value = 1575
$ws.Range("F33").Value = 1326
$ws.Range("F35").Value = 221
$ws.Range("I36").Value = "//i2.hdslb.com/bfs/openplatform/202404/A9SL4gE21713522395699.png"
$ws.Range("F37").Value = 244
$ws.Range("F39").Value = 4
$ws.Range("F41").Value = 1084
$ws.Range("F42").Value = 379
$ws.Range("F43").Value = 3365
$ws.Range("F44").Value = 654
$ws.Range("F46").Value = 1490
$ws.Range("F47").Value = 43
$ws.Range("F49").Value = 1
